$wb = $excel.ActiveWorkbook

# --- Sheet2 ("Sheet2 - Numbers"): extend the data range with a new column AA
# containing 30 values (100..129), one per existing row. This mirrors the
# companion "reader/sheet2" fixture referenced in the commit message. ---
$ws2 = $wb.Worksheets.Item("Sheet2 - Numbers")

for ($i = 1; $i -le 30; $i++) {
    $ws2.Range("AA" + $i).Value = 99 + $i
}

# Make Sheet2 the active sheet/tab and select the newly written column,
# matching the saved selection (AA1:AA30, active cell AA1).
$ws2.Activate()
$ws2.Range("AA1:AA30").Select()

# --- Sheet4 ("Sheet4 - Dates"): page setup now targets Letter-size paper. ---
$ws4 = $wb.Worksheets.Item("Sheet4 - Dates")
$ws4.PageSetup.PaperSize = 9
